$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.398.33"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3
$ws.Range("D3").Value = "1.848.35"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6303"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2957"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.57%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "

# Row 12
$ws.Range("D12").Value = "1.867.05"
$ws.Range("E12").Value = "  -1.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.993"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6850"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001004"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.147"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "

# Row 18
$ws.Range("D18").Value = "29.439.34"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.552"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "

# Row 23
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "

# Row 25
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1398"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.379"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "

# Row 28
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.467"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "

# Row 29
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05696"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.48%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.255"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

# Row 33
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.98%  "

# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.157"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.588"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.256.28"
$ws.Range("E37").Value = "  +1.42%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01817"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.12%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.785"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9119"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.194"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.0000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "2.030.95"
$ws.Range("E43").Value = "  -1.60%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.62%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.034"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.29%  "

# Row 47
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4024"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000117"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.106"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.696"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1129"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.14%  "
